# Apply the "Updated symbol list" data refresh (GitHub Actions run, 2023-01-12 08:18:53 UTC)
# Values are stored as text in the sheet (inline strings), so numeric-looking values
# (prices, volume percentages, hour) must be forced to Text format before assignment
# to avoid Excel auto-converting them to numbers/percentages and losing exact formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "284.33"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.42%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "8"
# Row 3
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.59%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "8"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.072"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.13%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "8"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06469"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.80%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "8"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.225"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.85%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "8"
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.321"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "12.01%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "8"
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9116"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.92%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "8"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1548"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.65%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "8"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06399"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "24.38%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "8"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07574"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.99%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "8"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02986"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.25%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "8"
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08947"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.33%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "8"
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001597"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.36%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "8"
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006398"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.69%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "8"
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006111"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.23%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "8"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.456"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.70%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "8"
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.370"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.84%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "8"
# Row 19
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "8"
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3189"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.27%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "8"
# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.45%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "8"
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.980"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.00%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "8"
# Row 23
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "8"
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04459"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.85%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "8"
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001184"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.56%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "8"
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004315"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "11.56%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "8"
# Row 27
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "8"
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001201"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-7.62%"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "8"
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0001638"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "-15.61%"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "8"
# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "8"
# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "8"
# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "8"
# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "8"
# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "8"
# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "8"
# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "8"
# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "8"
# Row 38
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "8"
# Row 39
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "8"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04135"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.65%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "8"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006738"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.49%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "8"
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1229"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.73%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "8"
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002101"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.55%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "8"
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01178"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.60%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "8"
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005412"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.94%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "8"
# Row 46
$ws.Range("B46").Value = "CoinbaseStockToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.01851"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.08%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "8"
# Row 47
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.041"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "21.16%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "8"
# Row 48
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "8"
# Row 49
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "8"
# Row 50
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "8"
# Row 51
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "8"
